$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10, shifting rows 10-13 down to 11-14
$ws.Rows.Item(10).Insert()

# Restore the cell borders for the new row (column A keeps the left/right
# "box" border, column E keeps the right border), matching the formatting
# used by the other data rows in the table.
$ws.Cells.Item(11, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(11, 5).Copy()
$ws.Cells.Item(10, 5).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Fill the new row 10 with the Taulman T-Glase profile data
# (values are set in the order that reproduces the shared-string table order)
$ws.Cells.Item(10, 1).Value = "Taulman T-Glase"
$ws.Cells.Item(10, 2).Value = "Taulman T-Glase"
$ws.Cells.Item(10, 4).Value = "Taulman T-Glase"
$ws.Cells.Item(10, 5).Value = "Standard Slic3r profile for T-Glase is nonsense, bed temperature above the glass transition for the material? ¯\_(ツ)_/¯ Generated new profile using guidance from Taulman, and a bit of trial and error"
$ws.Cells.Item(10, 3).Value = "Taulman T-Glase mod"

# Column E grows to fit the new (longer) note text
$ws.Columns.Item(5).ColumnWidth = 180.6

# Update the selection to match the post-edit workbook state
$ws.Range("C11").Select()
